$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - entry order matters because it controls the order in which
# strings land in the shared-strings table (Fecha, Horas, Requerimiento, ...)
$ws.Range("B1").Value = "Fecha"
$ws.Range("D1").Value = "Horas"
$ws.Range("C1").Value = "Requerimiento"

# Data row - set number format before value so the engine reuses the
# built-in date format (numFmtId 14) instead of inventing a custom one.
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = "2019-02-20"
$ws.Range("C2").Value = "Desarrollo carga de ingresos"

# Column widths (characters). The target XML widths are 20.28515625 and
# 97.42578125; the runtime stores ColumnWidth with its own padding/rounding
# (xml = round((cw + 5/6) * 6) / 6), so feed it the inverse to land on the
# closest representable grid point.
$ws.Columns.Item(2).ColumnWidth = 19.451822916666668
$ws.Columns.Item(3).ColumnWidth = 96.59244791666667

# Row height for the data row
$ws.Rows.Item(2).RowHeight = 30.75

# Selection left on the sheet after editing
$ws.Range("C6").Select()
